$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.863.81'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.304.59'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.08'
$ws.Range("E5").Value = '  +2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.04'
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -1.25%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.26'
$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0785'
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.66'
$ws.Range("E12").Value = '  +5.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.86'
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.662.54'
$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.294.92'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.780'
$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.762.76'
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("E19").Value = '  +1.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0895'
$ws.Range("E20").Value = '  -0.94%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.19'
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.04'
$ws.Range("E23").Value = '  -1.88%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  +0.98%  '

$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.72'
$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '166.80'
$ws.Range("E28").Value = '  +0.66%  '

$ws.Range("E29").Value = '  +1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.04'
$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.18'
$ws.Range("E31").Value = '  +1.11%  '

$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.05'
$ws.Range("E33").Value = '  +5.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.98'
$ws.Range("E34").Value = '  -0.29%  '

$ws.Range("E35").Value = '  -6.51%  '

$ws.Range("E36").Value = '  -0.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0686'
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("E41").Value = '  -1.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.999.14'
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0280'
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.23'
$ws.Range("E44").Value = '  +1.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.13'
$ws.Range("E45").Value = '  +6.35%  '

$ws.Range("E46").Value = '  +1.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.77'
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.55'
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.530.99'
$ws.Range("E49").Value = '  +0.63%  '

$ws.Range("E50").Value = '  +1.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.11'
$ws.Range("E51").Value = '  -0.86%  '
